# ============================================================
# Weekly crime-stat refresh (cs-en-us-020pct) for 020 Precinct
# - bumps the report Volume/Number and week-covering dates
# - refreshes the weekly/28-day/YTD/2-year crime figures (rows 15-30)
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Volume/Number and the week-covering date range ---
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# --- Crime-count cells whose style/number-format is unchanged ---
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = -28.571428571428
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 400
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 95
$ws.Range("J16").Value = 112
$ws.Range("K16").Value = -15.178571428571
$ws.Range("L16").Value = -12.844036697247
$ws.Range("M16").Value = 13.095238095238
$ws.Range("N16").Value = -83.362521891418
$ws.Range("C17").Value = 3
$ws.Range("F17").Value = 8
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 104
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 4
$ws.Range("L17").Value = 25.301204819277
$ws.Range("M17").Value = 96.226415094339
$ws.Range("N17").Value = 10.63829787234
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 105
$ws.Range("J18").Value = 134
$ws.Range("K18").Value = -21.641791044776
$ws.Range("L18").Value = 32.911392405063
$ws.Range("M18").Value = 6.060606060606
$ws.Range("N18").Value = -88.027366020524
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 33.333333333333
$ws.Range("I19").Value = 720
$ws.Range("J19").Value = 788
$ws.Range("K19").Value = -8.629441624365
$ws.Range("L19").Value = 0.558659217877
$ws.Range("M19").Value = 9.923664122137
$ws.Range("N19").Value = -58.549222797927
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 91
$ws.Range("K20").Value = 30
$ws.Range("L20").Value = 5.813953488372
$ws.Range("M20").Value = 175.757575757576
$ws.Range("N20").Value = -91.925465838509
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 5.555555555555
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = 25.974025974026
$ws.Range("I21").Value = 1129
$ws.Range("J21").Value = 1218
$ws.Range("K21").Value = -7.307060755336
$ws.Range("L21").Value = 4.537037037037
$ws.Range("M21").Value = 20.619658119658
$ws.Range("N21").Value = -74.497402304043
$ws.Range("I22").Value = 27
$ws.Range("J22").Value = 26
$ws.Range("K22").Value = 3.846153846153
$ws.Range("L22").Value = 17.391304347826
$ws.Range("M22").Value = -15.625
$ws.Range("C23").Value = 4
$ws.Range("F23").Value = 6
$ws.Range("I23").Value = 40
$ws.Range("K23").Value = 66.666666666666
$ws.Range("L23").Value = 66.666666666666
$ws.Range("M23").Value = 60
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 39.130434782608
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 41.891891891891
$ws.Range("I24").Value = 1148
$ws.Range("J24").Value = 1193
$ws.Range("K24").Value = -3.772003352891
$ws.Range("L24").Value = -5.592105263157
$ws.Range("M24").Value = 13.663366336633
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -28.571428571428
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -11.111111111111
$ws.Range("I25").Value = 217
$ws.Range("J25").Value = 219
$ws.Range("K25").Value = -0.913242009132
$ws.Range("L25").Value = 16.042780748663
$ws.Range("M25").Value = -20.220588235294
$ws.Range("D26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = 5.263157894736
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 41
$ws.Range("J27").Value = 53
$ws.Range("K27").Value = -22.641509433962
$ws.Range("L27").Value = -21.153846153846
$ws.Range("I30").Value = 16
$ws.Range("K30").Value = 300
$ws.Range("L30").Value = 220

# --- Cells switching from the blank-placeholder text ("0" / "***.*") to a real number ---
$ws.Range("D17").Value = 2
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("E17").Value = 50
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C20").Value = 1
$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'

# --- Cells switching from a real number back to the blank-placeholder text ---
# (NumberFormat "@" forces text storage for the assignment, then the format is
#  copied from an untouched donor cell that already carries the target General style
#  so the saved style index matches the original placeholder cells exactly.)
$donorZero = $ws.Range("D23")
$donorStar = $ws.Range("E23")
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$donorZero.Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$donorStar.Copy()
$ws.Range("E20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
